$d = $word.ActiveDocument

# Locate the last paragraph (current end-of-document content about 2024.8.28)
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$paraStart = $lastPara.Range.Start
$paraEnd = $lastPara.Range.End

# Remove the whole last paragraph (text + its paragraph mark)
$delRange = $d.Range($paraStart, $paraEnd)
$delRange.Delete()

# Insert the corrected version of that paragraph (merged runs, rFonts hint fixed)
$insertPoint = $d.Range($paraStart, $paraStart)
$para38Xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:tabs><w:tab w:val="left" w:pos="2107"/></w:tabs><w:bidi w:val="0"/><w:ind w:firstLine="420" w:firstLineChars="200"/><w:jc w:val="left"/><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="宋体" w:hAnsi="宋体" w:eastAsia="宋体" w:cs="宋体"/><w:b w:val="0"/><w:bCs w:val="0"/><w:kern w:val="2"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="宋体" w:hAnsi="宋体" w:eastAsia="宋体" w:cs="宋体"/><w:b w:val="0"/><w:bCs w:val="0"/><w:kern w:val="2"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>今天我们继续小程序的完善工作。我们为了避免昨天的错误，决定重新编写菜谱页面的代码。在选择图片时我们通过学习github中的教学视频实现了通过文件夹引入照片，而不是利用组件通过复制图片链接引入照片。同时为防止图片过大，我选择将图片图像大小进行压缩，使其在10kb左右，以此避免图片大小超过200kb的要求。此外，在实现跳转界面时，我更改了以前组件按钮的跳转的方法，而是通过bind:tap指令进行跳转。在完成代码编写后，我尝试进行上传发现可以正常上传，并用手机检查了菜谱页面的跳转功能及其能否正常显示菜谱。发现功能正常，可以实现。对于之后，因为选择食堂看菜谱的页面过于简单，我计划继续进行菜谱页面的美观优化以及排版设计，使其更加美观耐看、引人注目。</w:t></w:r></w:p>
"@
$insertPoint.InsertXML($para38Xml)

# Insert the new diary entry: 4 blank paragraphs + date heading + new content paragraph
$end2 = $d.Content.End
$insertPoint2 = $d.Range($end2 - 1, $end2 - 1)
$newBlockXml = @"
<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="2107"/></w:tabs><w:bidi w:val="0"/><w:ind w:firstLine="420" w:firstLineChars="200"/><w:jc w:val="left"/><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="宋体" w:hAnsi="宋体" w:eastAsia="宋体" w:cs="宋体"/><w:b w:val="0"/><w:bCs w:val="0"/><w:kern w:val="2"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="2107"/></w:tabs><w:bidi w:val="0"/><w:ind w:firstLine="420" w:firstLineChars="200"/><w:jc w:val="left"/><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="宋体" w:hAnsi="宋体" w:eastAsia="宋体" w:cs="宋体"/><w:b w:val="0"/><w:bCs w:val="0"/><w:kern w:val="2"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="2107"/></w:tabs><w:bidi w:val="0"/><w:ind w:firstLine="420" w:firstLineChars="200"/><w:jc w:val="left"/><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="宋体" w:hAnsi="宋体" w:eastAsia="宋体" w:cs="宋体"/><w:b w:val="0"/><w:bCs w:val="0"/><w:kern w:val="2"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="2107"/></w:tabs><w:bidi w:val="0"/><w:ind w:firstLine="420" w:firstLineChars="200"/><w:jc w:val="left"/><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="宋体" w:hAnsi="宋体" w:eastAsia="宋体" w:cs="宋体"/><w:b w:val="0"/><w:bCs w:val="0"/><w:kern w:val="2"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="2107"/></w:tabs><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="宋体" w:hAnsi="宋体" w:eastAsia="宋体" w:cs="宋体"/><w:b/><w:bCs/><w:kern w:val="2"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="宋体" w:hAnsi="宋体" w:eastAsia="宋体" w:cs="宋体"/><w:b/><w:bCs/><w:kern w:val="2"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>2024.8.29</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="宋体" w:hAnsi="宋体" w:eastAsia="宋体" w:cs="宋体"/><w:b/><w:bCs/><w:kern w:val="2"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t xml:space="preserve">   天气晴</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="2107"/></w:tabs><w:bidi w:val="0"/><w:ind w:firstLine="420" w:firstLineChars="200"/><w:jc w:val="left"/><w:rPr><w:rFonts w:hint="default" w:ascii="宋体" w:hAnsi="宋体" w:eastAsia="宋体" w:cs="宋体"/><w:b w:val="0"/><w:bCs w:val="0"/><w:kern w:val="2"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="宋体" w:hAnsi="宋体" w:eastAsia="宋体" w:cs="宋体"/><w:b w:val="0"/><w:bCs w:val="0"/><w:kern w:val="2"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>今天我们继续小程序的完善工作。根据昨天的计划，我们进行了菜谱页面的美观优化以及排版设计。通过加入图片、改变背景颜色等方式使小程序页面更加漂亮、引人注目。此外，我们还进行了小程序功能的扩充，通过学习github中有关订单的代码知识，我们增加了代购功能，用户可以通过小程序对食堂菜品进行线上购买并送货上门。最后，我们将完善好的随机选择功能、显示菜谱功能、用户反馈功能以及代购服务功能的代码进行合并，完成了小程序的编写。并且成功上传，手机可以查看小程序，发现功能正常。明天，我们将进行对小程序的测试与完善同时完成答辩要求的ppt。</w:t></w:r></w:p>
"@
$insertPoint2.InsertXML($newBlockXml)

Write-Output ("Done. ParaCount=" + $d.Paragraphs.Count)
